$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 72, pushing the existing row 72 (and below) down.
$ws.Rows.Item(72).Insert()

# Populate the newly inserted row 72 with this week's price record.
$ws.Cells.Item(72, 1).Value() = 10
$ws.Cells.Item(72, 2).Value() = "Vega Modelo de Temuco"
$ws.Cells.Item(72, 3).Value() = "La Araucanía"
$ws.Cells.Item(72, 4).Value() = 44448
$ws.Cells.Item(72, 5).Value() = 9
$ws.Cells.Item(72, 6).Value() = 100112012
$ws.Cells.Item(72, 7).Value() = "Espinaca"
$ws.Cells.Item(72, 8).Value() = "Sin especificar"
$ws.Cells.Item(72, 9).Value() = "Primera"
$ws.Cells.Item(72, 10).Value() = 75
$ws.Cells.Item(72, 11).Value() = 9000
$ws.Cells.Item(72, 12).Value() = 9000
$ws.Cells.Item(72, 13).Value() = 9000
$ws.Cells.Item(72, 14).Value() = "`$/docena de atados"
$ws.Cells.Item(72, 15).Value() = "Región de La Araucanía"
$ws.Cells.Item(72, 16).Value() = 3000
$ws.Cells.Item(72, 17).Value() = 3
$ws.Cells.Item(72, 18).Value() = "Hortaliza"
